$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 53, shifting existing rows 53:58 down to 54:59.
$ws.Rows("53:53").Insert(-4121)

# New row 53 keeps the same descriptive columns as the row below it (old row 53,
# now shifted to row 54), only the measurement columns change.
$ws.Range("A53").Value = 9
$ws.Range("B53").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C53").Value = "Metropolitana"
$ws.Range("D53").Value = 44505
$ws.Range("D53").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E53").Value = 13
$ws.Range("F53").Value = "Fruta"
$ws.Range("G53").Value = 100102
$ws.Range("H53").Value = "Cítricos"
$ws.Range("I53").Value = 100102006
$ws.Range("J53").Value = "Pomelo"
$ws.Range("K53").Value = "Start Ruby"
$ws.Range("L53").Value = "Primera"
$ws.Range("M53").Value = 500
$ws.Range("N53").Value = 7500
$ws.Range("O53").Value = 8000
$ws.Range("P53").Value = 7720
$ws.Range("Q53").Value = "$/caja 14 kilos granel"
$ws.Range("R53").Value = "Región Metropolitana"
$ws.Range("S53").Value = 551
$ws.Range("T53").Value = 14
